$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was empty) -> Alvearie Team
$ws.Range("B9").Value = "Alvearie Team"

# Row 10: Contact / No display for ContactDetail -> Jurisdiction / United States of America
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the old duplicate "Contact" row (row 11); rows below shift up by one
$ws.Rows.Item(11).Delete()

# Case Sensitive value (now row 14 after the row shift) -> "true" (stored as text, not boolean)
$ws.Range("Z1").Formula = "=""true"""
$ws.Range("Z1").Copy()
$ws.Range("B14").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
